$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as Text (matches source inlineStr cells)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.665.96'
$ws.Range("E2").Value = '  +2.75%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.567.32'
$ws.Range("E3").Value = '  +2.57%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.21'
$ws.Range("E5").Value = '  +2.85%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '186.73'
$ws.Range("E6").Value = '  +2.57%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.630'
$ws.Range("E7").Value = '  +2.92%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.556.72'
$ws.Range("E8").Value = '  +2.42%  '

$ws.Range("E9").Value = '  -0.03%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.221'
$ws.Range("E10").Value = '  +20.83%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.653'
$ws.Range("E11").Value = '  +2.08%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.64'
$ws.Range("E12").Value = '  +1.93%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000318'
$ws.Range("E13").Value = '  +6.52%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.49'
$ws.Range("E14").Value = '  +1.17%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.131.74'
$ws.Range("E15").Value = '  +2.39%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '70.651.36'
$ws.Range("E16").Value = '  +3.02%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.23'
$ws.Range("E17").Value = '  +0.64%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.537.76'
$ws.Range("E18").Value = '  +1.96%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.68'
$ws.Range("E19").Value = '  +3.63%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '574.07'
$ws.Range("E20").Value = '  +6.84%  '

$ws.Range("E22").Value = '  -0.24%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.65'
$ws.Range("E23").Value = '  -7.83%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.55'
$ws.Range("E24").Value = '  +4.41%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.89'
$ws.Range("E25").Value = '  -1.26%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '94.57'
$ws.Range("E26").Value = '  +0.66%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.31'
$ws.Range("E27").Value = '  +5.40%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.96'
$ws.Range("E28").Value = '  +2.62%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.23'
$ws.Range("E29").Value = '  +2.71%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.52'
$ws.Range("E30").Value = '  +3.80%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.21'
$ws.Range("E31").Value = '  +0.88%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.31'
$ws.Range("E32").Value = '  -1.28%  '

$ws.Range("E33").Value = '  +3.62%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.44'
$ws.Range("E34").Value = '  +14.71%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '63.12'
$ws.Range("E35").Value = '  -1.87%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '549.92'
$ws.Range("E36").Value = '  -3.20%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.416'
$ws.Range("E37").Value = '  +5.74%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.43'
$ws.Range("E38").Value = '  +10.78%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.97'
$ws.Range("E39").Value = '  +0.98%  '

$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0802'
$ws.Range("E40").Value = '  +5.60%  '

$ws.Range("B41").Value = 'Dai'
$ws.Range("C41").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.600.77'
$ws.Range("E42").Value = '  +12.16%  '

$ws.Range("E43").Value = '  +3.85%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.43'
$ws.Range("E44").Value = '  +3.63%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0469'
$ws.Range("E45").Value = '  +8.16%  '

$ws.Range("E46").Value = '  +0.92%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.92'
$ws.Range("E47").Value = '  -0.89%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.34'
$ws.Range("E48").Value = '  +4.22%  '

$ws.Range("E49").Value = '  +3.19%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.54'
$ws.Range("E50").Value = '  +17.86%  '

$ws.Range("B51").Value = 'FLOKI'
$ws.Range("C51").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000264'
$ws.Range("E51").Value = '  +17.11%  '
